# Add a new "UK" Test Data sheet, cloned from the existing "Poland" sheet,
# populate its market-specific values, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Poland")

# Copy the Poland sheet and place the copy immediately after it.
$source.Copy([System.Type]::Missing, $source) | Out-Null

# Excel names the freshly copied sheet "Poland (2)" - rename it to "UK".
$newSheet = $wb.Worksheets.Item("Poland (2)")
$newSheet.Name = "UK"

# Fill in the UK-specific values (code first, then market name, to match
# the order new shared strings are appended to the workbook).
$newSheet.Range("B4").Value = "NGC-2741/T3334"
$newSheet.Range("B2").Value = "UK Market"

# Make the new sheet the active tab with B4 selected, like the source diff.
$newSheet.Activate() | Out-Null
$newSheet.Range("B4").Select() | Out-Null
